$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041975869641576
$ws.Cells.Item(2, 4).Value = 1.049034107013959
$ws.Cells.Item(2, 5).Value = 1.040167714461241
$ws.Cells.Item(2, 6).Value = 1.057516302088566
$ws.Cells.Item(2, 9).Value = 1.039430768745064
$ws.Cells.Item(2, 10).Value = 1.047054180259533
$ws.Cells.Item(2, 11).Value = 1.051792290542114
$ws.Cells.Item(2, 12).Value = 1.042950830137768
$ws.Cells.Item(2, 13).Value = 1.060251057777758
$ws.Cells.Item(2, 14).Value = 1.019541078125463
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.043135185872939
$ws.Cells.Item(3, 4).Value = 1.050081092895724
$ws.Cells.Item(3, 5).Value = 1.041158675696671
$ws.Cells.Item(3, 6).Value = 1.058655087041409
$ws.Cells.Item(3, 9).Value = 1.039678457180744
$ws.Cells.Item(3, 10).Value = 1.047858701676387
$ws.Cells.Item(3, 11).Value = 1.05265083460236
$ws.Cells.Item(3, 12).Value = 1.043751684027528
$ws.Cells.Item(3, 13).Value = 1.061202869151823
$ws.Cells.Item(3, 14).Value = 1.019810979323055
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043884949076658
$ws.Cells.Item(4, 4).Value = 1.050758476351873
$ws.Cells.Item(4, 5).Value = 1.041799912646078
$ws.Cells.Item(4, 6).Value = 1.0593918905822
$ws.Cells.Item(4, 9).Value = 1.039836649178255
$ws.Cells.Item(4, 10).Value = 1.048378356864655
$ws.Cells.Item(4, 11).Value = 1.053205682514951
$ws.Cells.Item(4, 12).Value = 1.044269297905248
$ws.Cells.Item(4, 13).Value = 1.061818104173962
$ws.Cells.Item(4, 14).Value = 1.019985225648882
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.044200057025434
$ws.Cells.Item(5, 4).Value = 1.051043228251818
$ws.Cells.Item(5, 5).Value = 1.042069493580807
$ws.Cells.Item(5, 6).Value = 1.059701627607325
$ws.Cells.Item(5, 9).Value = 1.039902655405306
$ws.Cells.Item(5, 10).Value = 1.048596599257871
$ws.Cells.Item(5, 11).Value = 1.053438776537513
$ws.Cells.Item(5, 12).Value = 1.044486761476266
$ws.Cells.Item(5, 13).Value = 1.06207659383219
$ws.Cells.Item(5, 14).Value = 1.020058383692768
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.044252959650586
$ws.Cells.Item(6, 4).Value = 1.051091038186652
$ws.Cells.Item(6, 5).Value = 1.042114757706865
$ws.Cells.Item(6, 6).Value = 1.059753632947464
$ws.Cells.Item(6, 9).Value = 1.039913708957117
$ws.Cells.Item(6, 10).Value = 1.048633230192977
$ws.Cells.Item(6, 11).Value = 1.053477904446027
$ws.Cells.Item(6, 12).Value = 1.044523266293777
$ws.Cells.Item(6, 13).Value = 1.062119986302884
$ws.Cells.Item(6, 14).Value = 1.020070661676359
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.043889159926318
$ws.Cells.Item(7, 4).Value = 1.050762281298982
$ws.Cells.Item(7, 5).Value = 1.041803514779865
$ws.Cells.Item(7, 6).Value = 1.059396029362254
$ws.Cells.Item(7, 9).Value = 1.039837533112255
$ws.Cells.Item(7, 10).Value = 1.04838127389652
$ws.Cells.Item(7, 11).Value = 1.053208797772585
$ws.Cells.Item(7, 12).Value = 1.044272204218798
$ws.Cells.Item(7, 13).Value = 1.061821558734304
$ws.Cells.Item(7, 14).Value = 1.019986203563822
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042367747938666
$ws.Cells.Item(8, 4).Value = 1.049387959042641
$ws.Cells.Item(8, 5).Value = 1.040502610825296
$ws.Cells.Item(8, 6).Value = 1.057901174001034
$ws.Cells.Item(8, 9).Value = 1.039514906695736
$ws.Cells.Item(8, 10).Value = 1.047326263686803
$ws.Cells.Item(8, 11).Value = 1.052082582212811
$ws.Cells.Item(8, 12).Value = 1.043221605343265
$ws.Cells.Item(8, 13).Value = 1.060572861963024
$ws.Cells.Item(8, 14).Value = 1.019632374939599
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039683763980417
$ws.Cells.Item(9, 4).Value = 1.046965538709641
$ws.Cells.Item(9, 5).Value = 1.03821037550059
$ws.Cells.Item(9, 6).Value = 1.055266505129249
$ws.Cells.Item(9, 9).Value = 1.038930473396578
$ws.Cells.Item(9, 10).Value = 1.045460098900504
$ws.Cells.Item(9, 11).Value = 1.050092755578374
$ws.Cells.Item(9, 12).Value = 1.041365759568871
$ws.Cells.Item(9, 13).Value = 1.058367487972369
$ws.Cells.Item(9, 14).Value = 1.019005831285246
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037892295753096
$ws.Cells.Item(10, 4).Value = 1.045350087305054
$ws.Cells.Item(10, 5).Value = 1.03668227031211
$ws.Cells.Item(10, 6).Value = 1.053509643607718
$ws.Cells.Item(10, 9).Value = 1.038530139823522
$ws.Cells.Item(10, 10).Value = 1.044211175266681
$ws.Cells.Item(10, 11).Value = 1.048762608407838
$ws.Cells.Item(10, 12).Value = 1.040125432984886
$ws.Cells.Item(10, 13).Value = 1.056893827040334
$ws.Cells.Item(10, 14).Value = 1.018586074062078
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037116040955934
$ws.Cells.Item(11, 4).Value = 1.044650449076259
$ws.Cells.Item(11, 5).Value = 1.036020588118196
$ws.Cells.Item(11, 6).Value = 1.052748792181406
$ws.Cells.Item(11, 9).Value = 1.038354248614126
$ws.Cells.Item(11, 10).Value = 1.043669226190395
$ws.Cells.Item(11, 11).Value = 1.048185777454769
$ws.Cells.Item(11, 12).Value = 1.039587615351635
$ws.Cells.Item(11, 13).Value = 1.056254895839991
$ws.Cells.Item(11, 14).Value = 1.018403823205606
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036827622774995
$ws.Cells.Item(12, 4).Value = 1.044390550683542
$ws.Cells.Item(12, 5).Value = 1.035774808699028
$ws.Cells.Item(12, 6).Value = 1.052466159042067
$ws.Cells.Item(12, 9).Value = 1.038288532315192
$ws.Cells.Item(12, 10).Value = 1.043467747421451
$ws.Cells.Item(12, 11).Value = 1.047971385339609
$ws.Cells.Item(12, 12).Value = 1.039387732740934
$ws.Cells.Item(12, 13).Value = 1.056017443169771
$ws.Cells.Item(12, 14).Value = 1.018336052775041
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036889493221266
$ws.Cells.Item(13, 4).Value = 1.044446300764843
$ws.Cells.Item(13, 5).Value = 1.035827529285693
$ws.Cells.Item(13, 6).Value = 1.052526785690908
$ws.Cells.Item(13, 9).Value = 1.038302645996086
$ws.Cells.Item(13, 10).Value = 1.043510973224413
$ws.Cells.Item(13, 11).Value = 1.048017379138133
$ws.Cells.Item(13, 12).Value = 1.039430613374896
$ws.Cells.Item(13, 13).Value = 1.056068383254671
$ws.Cells.Item(13, 14).Value = 1.018350593130869
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037092201914735
$ws.Cells.Item(14, 4).Value = 1.044628966240581
$ws.Cells.Item(14, 5).Value = 1.036000271943726
$ws.Cells.Item(14, 6).Value = 1.052725430033357
$ws.Cells.Item(14, 9).Value = 1.038348824288119
$ws.Cells.Item(14, 10).Value = 1.043652575462859
$ws.Cells.Item(14, 11).Value = 1.048168058421967
$ws.Cells.Item(14, 12).Value = 1.039571095306555
$ws.Cells.Item(14, 13).Value = 1.05623527048772
$ws.Cells.Item(14, 14).Value = 1.018398222798795
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037217086431291
$ws.Cells.Item(15, 4).Value = 1.044741509576074
$ws.Cells.Item(15, 5).Value = 1.036106704214487
$ws.Cells.Item(15, 6).Value = 1.052847818809549
$ws.Cells.Item(15, 9).Value = 1.038377225569167
$ws.Cells.Item(15, 10).Value = 1.043739798090832
$ws.Cells.Item(15, 11).Value = 1.048260879459605
$ws.Cells.Item(15, 12).Value = 1.039657635837237
$ws.Cells.Item(15, 13).Value = 1.056338078602125
$ws.Cells.Item(15, 14).Value = 1.018427559149891
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037943801748895
$ws.Cells.Item(16, 4).Value = 1.045396517015889
$ws.Cells.Item(16, 5).Value = 1.036726183857554
$ws.Cells.Item(16, 6).Value = 1.053560136206456
$ws.Cells.Item(16, 9).Value = 1.038541759521132
$ws.Cells.Item(16, 10).Value = 1.044247118187384
$ws.Cells.Item(16, 11).Value = 1.048800872401098
$ws.Cells.Item(16, 12).Value = 1.040161110340616
$ws.Cells.Item(16, 13).Value = 1.056936213302881
$ws.Cells.Item(16, 14).Value = 1.018598159038694
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038399505916844
$ws.Cells.Item(17, 4).Value = 1.045807348464548
$ws.Cells.Item(17, 5).Value = 1.037114765819526
$ws.Cells.Item(17, 6).Value = 1.054006921446063
$ws.Cells.Item(17, 9).Value = 1.038644286015723
$ws.Cells.Item(17, 10).Value = 1.044565036285078
$ws.Cells.Item(17, 11).Value = 1.049139362537243
$ws.Cells.Item(17, 12).Value = 1.040476725754611
$ws.Cells.Item(17, 13).Value = 1.057311185450935
$ws.Cells.Item(17, 14).Value = 1.018705039621304
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038665258739171
$ws.Cells.Item(18, 4).Value = 1.046046966336627
$ws.Cells.Item(18, 5).Value = 1.037341418910956
$ws.Cells.Item(18, 6).Value = 1.054267512336244
$ws.Cells.Item(18, 9).Value = 1.038703842530606
$ws.Cells.Item(18, 10).Value = 1.044750360972168
$ws.Cells.Item(18, 11).Value = 1.049336714377638
$ws.Cells.Item(18, 12).Value = 1.040660746765959
$ws.Cells.Item(18, 13).Value = 1.057529820705567
$ws.Cells.Item(18, 14).Value = 1.018767333686226
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038755864870012
$ws.Cells.Item(19, 4).Value = 1.046128667636731
$ws.Cells.Item(19, 5).Value = 1.037418701749566
$ws.Cells.Item(19, 6).Value = 1.054356365201281
$ws.Cells.Item(19, 9).Value = 1.038724108137057
$ws.Cells.Item(19, 10).Value = 1.044813533007098
$ws.Cells.Item(19, 11).Value = 1.049403992081729
$ws.Cells.Item(19, 12).Value = 1.040723480967062
$ws.Cells.Item(19, 13).Value = 1.057604356197692
$ws.Cells.Item(19, 14).Value = 1.018788566291252
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038350618534097
$ws.Cells.Item(20, 4).Value = 1.045763271506522
$ws.Cells.Item(20, 5).Value = 1.037073074678054
$ws.Cells.Item(20, 6).Value = 1.053958986803761
$ws.Cells.Item(20, 9).Value = 1.038633311283354
$ws.Cells.Item(20, 10).Value = 1.044530938227754
$ws.Cells.Item(20, 11).Value = 1.049103054416161
$ws.Cells.Item(20, 12).Value = 1.040442870673003
$ws.Cells.Item(20, 13).Value = 1.057270962734497
$ws.Cells.Item(20, 14).Value = 1.018693577268481
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037032511551776
$ws.Cells.Item(21, 4).Value = 1.044575176415259
$ws.Cells.Item(21, 5).Value = 1.035949403575828
$ws.Cells.Item(21, 6).Value = 1.052666934760911
$ws.Cells.Item(21, 9).Value = 1.038335236492192
$ws.Cells.Item(21, 10).Value = 1.043610881964582
$ws.Cells.Item(21, 11).Value = 1.048123690763783
$ws.Cells.Item(21, 12).Value = 1.039529730010406
$ws.Cells.Item(21, 13).Value = 1.056186129820905
$ws.Cells.Item(21, 14).Value = 1.018384199103104
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036203286621689
$ws.Cells.Item(22, 4).Value = 1.043828048382431
$ws.Cells.Item(22, 5).Value = 1.035242899659118
$ws.Cells.Item(22, 6).Value = 1.051854459585335
$ws.Cells.Item(22, 9).Value = 1.038145612047236
$ws.Cells.Item(22, 10).Value = 1.043031394542792
$ws.Cells.Item(22, 11).Value = 1.0475071651607
$ws.Cells.Item(22, 12).Value = 1.038954946911861
$ws.Cells.Item(22, 13).Value = 1.055503327574019
$ws.Cells.Item(22, 14).Value = 1.01818925049279
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036642920455376
$ws.Cells.Item(23, 4).Value = 1.044224127322182
$ws.Cells.Item(23, 5).Value = 1.035617431716226
$ws.Cells.Item(23, 6).Value = 1.052285178950798
$ws.Cells.Item(23, 9).Value = 1.038246345403508
$ws.Cells.Item(23, 10).Value = 1.04333868798764
$ws.Cells.Item(23, 11).Value = 1.047834069480236
$ws.Cells.Item(23, 12).Value = 1.039259712697889
$ws.Cells.Item(23, 13).Value = 1.05586536300873
$ws.Cells.Item(23, 14).Value = 1.018292637309022
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038372708797259
$ws.Cells.Item(24, 4).Value = 1.04578318802689
$ws.Cells.Item(24, 5).Value = 1.037091913110595
$ws.Cells.Item(24, 6).Value = 1.053980646440051
$ws.Cells.Item(24, 9).Value = 1.038638271050776
$ws.Cells.Item(24, 10).Value = 1.044546346016492
$ws.Cells.Item(24, 11).Value = 1.049119460751584
$ws.Cells.Item(24, 12).Value = 1.040458168549643
$ws.Cells.Item(24, 13).Value = 1.05728913789444
$ws.Cells.Item(24, 14).Value = 1.018698756759214
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040378009922732
$ws.Cells.Item(25, 4).Value = 1.047591878287182
$ws.Cells.Item(25, 5).Value = 1.038802961823281
$ws.Cells.Item(25, 6).Value = 1.055947699680824
$ws.Cells.Item(25, 9).Value = 1.039083450986031
$ws.Cells.Item(25, 10).Value = 1.045943392809252
$ws.Cells.Item(25, 11).Value = 1.050607804414758
$ws.Cells.Item(25, 12).Value = 1.041846083461481
$ws.Cells.Item(25, 13).Value = 1.058938228161684
$ws.Cells.Item(25, 14).Value = 1.019168170497726

Write-Output "Updated 264 cells (B:F, I:N) for rows 2-25"
